$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R (2022) into the new column S (2023) for the
# header/data rows of the table, then fill in the 2023 figures.
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("S3").Value = 2023

$ws.Range("S4").Value = 58.6
$ws.Range("S5").Value = 58.6
$ws.Range("S6").Value = 1294
$ws.Range("S7").Value = 1057
$ws.Range("S8").Value = 976.4
$ws.Range("S9").Value = 35.2
$ws.Range("S10").Value = 12.4
$ws.Range("S11").Value = 23.4
$ws.Range("S12").Value = 2
$ws.Range("S13").Value = 35.2
$ws.Range("S14").Value = "_"

# Match the saved selection from the source workbook.
$ws.Range("G21").Select()
